# Scheduled-runner price refresh: updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns on each Leve-profit worksheet with freshly pulled market-board data.
$wb = $excel.ActiveWorkbook

# ALC!row15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 554.29785
$ws.Range("I15").Value = 554.29785
$ws.Range("K15").Value = 1662.89355
$ws.Range("M15").Value = -1493.89355

# ALC!row70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 74878.37
$ws.Range("I70").Value = 1050
$ws.Range("J70").Value = 117066
$ws.Range("K70").Value = 3150
$ws.Range("L70").Value = 351198
$ws.Range("M70").Value = -2880
$ws.Range("N70").Value = -351738

# ALC!row73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 74878.37
$ws.Range("I73").Value = 1050
$ws.Range("J73").Value = 117066
$ws.Range("K73").Value = 3150
$ws.Range("L73").Value = 351198
$ws.Range("M73").Value = -2214
$ws.Range("N73").Value = -353070

# ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6432.7896
$ws.Range("I76").Value = 5684.4287
$ws.Range("J76").Value = 6869.3335
$ws.Range("K76").Value = 5684.4287
$ws.Range("L76").Value = 6869.3335
$ws.Range("M76").Value = -5369.4287
$ws.Range("N76").Value = -7499.3335

# ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6432.7896
$ws.Range("I79").Value = 5684.4287
$ws.Range("J79").Value = 6869.3335
$ws.Range("K79").Value = 5684.4287
$ws.Range("L79").Value = 6869.3335
$ws.Range("M79").Value = -4592.4287
$ws.Range("N79").Value = -9053.333500000001

# ALC!row106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 20219.334
$ws.Range("I106").Value = 22463.2
$ws.Range("K106").Value = 22463.2
$ws.Range("M106").Value = -21832.2

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5409.294
$ws.Range("J138").Value = 5574.0215
$ws.Range("L138").Value = 16722.0645
$ws.Range("N138").Value = -27002.0645

# ALC!row141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2015
$ws.Range("I141").Value = 2015
$ws.Range("K141").Value = 6045
$ws.Range("M141").Value = -865

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1951.3
$ws.Range("I45").Value = 1923.6666
$ws.Range("K45").Value = 1923.6666
$ws.Range("M45").Value = -1546.6666

# ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6718.769
$ws.Range("J63").Value = 7168.125
$ws.Range("L63").Value = 7168.125
$ws.Range("N63").Value = -8540.125

# ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6718.769
$ws.Range("J66").Value = 7168.125
$ws.Range("L66").Value = 35840.625
$ws.Range("N66").Value = -42704.625

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2135.4285
$ws.Range("I102").Value = 1066.3334
$ws.Range("J102").Value = 2937.25
$ws.Range("K102").Value = 1066.3334
$ws.Range("L102").Value = 2937.25
$ws.Range("M102").Value = 555.6666
$ws.Range("N102").Value = -6181.25

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 558671.9399999999
$ws.Range("I122").Value = 911009
$ws.Range("J122").Value = 4999.4287
$ws.Range("K122").Value = 2733027
$ws.Range("L122").Value = 14998.2861
$ws.Range("M122").Value = -2730577
$ws.Range("N122").Value = -19898.2861

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4034.8064
$ws.Range("I105").Value = 3279
$ws.Range("K105").Value = 3279
$ws.Range("M105").Value = -1532

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1201.6666
$ws.Range("I134").Value = 681.8
$ws.Range("K134").Value = 2045.4
$ws.Range("M134").Value = 489.6000000000001

# CRP!row22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 322
$ws.Range("I22").Value = 299.66666
$ws.Range("K22").Value = 299.66666
$ws.Range("M22").Value = 50.33334000000002

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5576.4375
$ws.Range("I31").Value = 3366
$ws.Range("K31").Value = 3366
$ws.Range("M31").Value = -3071

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5576.4375
$ws.Range("I34").Value = 3366
$ws.Range("K34").Value = 3366
$ws.Range("M34").Value = -3164

# CRP!row94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1597
$ws.Range("J94").Value = 1597
$ws.Range("L94").Value = 1597
$ws.Range("N94").Value = -2499

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 6339.769
$ws.Range("I122").Value = 5853.727
$ws.Range("K122").Value = 17561.181
$ws.Range("M122").Value = -15111.181

# CRP!row125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 86500
$ws.Range("J125").Value = 86500
$ws.Range("L125").Value = 86500
$ws.Range("N125").Value = -91420

# CUL!row17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 481
$ws.Range("I17").Value = 189.5
$ws.Range("J17").Value = 918.25
$ws.Range("K17").Value = 568.5
$ws.Range("L17").Value = 2754.75
$ws.Range("M17").Value = -399.5
$ws.Range("N17").Value = -3092.75

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 956.8919
$ws.Range("I113").Value = 1531.1666
$ws.Range("K113").Value = 4593.4998
$ws.Range("M113").Value = -2423.4998

# CUL!row117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2499.1667
$ws.Range("J117").Value = 4065
$ws.Range("L117").Value = 12195
$ws.Range("N117").Value = -19079

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4861.6113
$ws.Range("I131").Value = 2865.2
$ws.Range("J131").Value = 7357.125
$ws.Range("K131").Value = 8595.599999999999
$ws.Range("L131").Value = 22071.375
$ws.Range("M131").Value = -3555.599999999999
$ws.Range("N131").Value = -32151.375

# CUL!row134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1042.3334
$ws.Range("I134").Value = 1042.3334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3127.0002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 1942.9998
$ws.Range("N134").ClearContents()

# GSM!row97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1918.1154
$ws.Range("I97").Value = 1855.3334
$ws.Range("K97").Value = 1855.3334
$ws.Range("M97").Value = -1359.3334

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1192.6666
$ws.Range("I102").Value = 496.5
$ws.Range("J102").Value = 3778.4285
$ws.Range("K102").Value = 496.5
$ws.Range("L102").Value = 3778.4285
$ws.Range("M102").Value = 1125.5
$ws.Range("N102").Value = -7022.4285

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1003103.25
$ws.Range("I122").Value = 168524
$ws.Range("K122").Value = 505572
$ws.Range("M122").Value = -503122

# GSM!row136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 26402.643
$ws.Range("J136").Value = 26402.643
$ws.Range("L136").Value = 79207.929
$ws.Range("N136").Value = -84307.929

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1655.6111
$ws.Range("I7").Value = 1641.3529
$ws.Range("J7").Value = 1898
$ws.Range("K7").Value = 1641.3529
$ws.Range("L7").Value = 1898
$ws.Range("M7").Value = -1529.3529
$ws.Range("N7").Value = -2122

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 640
$ws.Range("I22").Value = 640
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 640
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -345
$ws.Range("N22").ClearContents()

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 640
$ws.Range("I27").Value = 640
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 640
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -533
$ws.Range("N27").ClearContents()

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2508.3333
$ws.Range("I46").Value = 825
$ws.Range("J46").Value = 3350
$ws.Range("K46").Value = 825
$ws.Range("L46").Value = 3350
$ws.Range("M46").Value = -637
$ws.Range("N46").Value = -3726

# LTW!row56
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 10051
$ws.Range("I56").Value = 10051
$ws.Range("K56").Value = 10051
$ws.Range("M56").Value = -9360

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2899.5881
$ws.Range("I68").Value = 2433.6667
$ws.Range("J68").Value = 2999.4285
$ws.Range("K68").Value = 2433.6667
$ws.Range("L68").Value = 2999.4285
$ws.Range("M68").Value = -1684.6667
$ws.Range("N68").Value = -4497.4285

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2899.5881
$ws.Range("I71").Value = 2433.6667
$ws.Range("J71").Value = 2999.4285
$ws.Range("K71").Value = 12168.3335
$ws.Range("L71").Value = 14997.1425
$ws.Range("M71").Value = -8424.333500000001
$ws.Range("N71").Value = -22485.1425

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1655.6111
$ws.Range("I126").Value = 1641.3529
$ws.Range("J126").Value = 1898
$ws.Range("K126").Value = 4924.0587
$ws.Range("L126").Value = 5694
$ws.Range("M126").Value = -2454.0587
$ws.Range("N126").Value = -10634

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2113.3
$ws.Range("I136").Value = 1517.25
$ws.Range("K136").Value = 4551.75
$ws.Range("M136").Value = -2001.75

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1356.8695
$ws.Range("I132").Value = 1391.0454
$ws.Range("K132").Value = 4173.1362
$ws.Range("M132").Value = -1643.1362

# WVR!row135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 79250
$ws.Range("J135").Value = 79250
$ws.Range("L135").Value = 79250
$ws.Range("N135").Value = -89390

